$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New header columns: Middle Name (Q) and Nationality (R)
$ws.Range("Q1").Value = "Middle Name"
$ws.Range("R1").Value = "Nationality"

# Row 2 data
$ws.Range("Q2").Value = "Sir"
$ws.Range("R2").Value = "India"

# Row 3 data
$ws.Range("Q3").Value = "Sir"
$ws.Range("R3").Value = "US"

# Update the view: scroll so column C is the left-most visible column,
# then select cell R4 (mirrors the author's recorded sheetView state).
$excel.Goto($ws.Range("C1"), $false)
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("R4").Select()
